$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# --- Cells changing from shared-string placeholder to numeric (copy style from stable donor, then set numeric value) ---
# Donor cells: J14 (style 14, numeric #,##0) and K14 (style 15, numeric #,##0.0 w/ negative paren)
$ws.Range("J14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("J14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("K14").Copy($ws.Range("N15"))
$ws.Range("N15").Value = 0
$ws.Range("J14").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$ws.Range("K14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 0
$ws.Range("J14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("J14").Copy($ws.Range("F27"))
$ws.Range("F27").Value = 1
$ws.Range("J14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("J14").Copy($ws.Range("G28"))
$ws.Range("G28").Value = 1
$ws.Range("K14").Copy($ws.Range("H28"))
$ws.Range("H28").Value = 0
$ws.Range("J14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("K14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("J14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100
$ws.Range("J14").Copy($ws.Range("J31"))
$ws.Range("J31").Value = 1
$ws.Range("K14").Copy($ws.Range("K31"))
$ws.Range("K31").Value = -100

# --- Cells changing from numeric to shared-string placeholder (copy style+text from stable donor) ---
# Donor cells: C14 (style 13, shared string "0") and E14 (style 13, shared string "***.*")
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("E14").Copy($ws.Range("H27"))
$ws.Range("C14").Copy($ws.Range("C28"))

# --- Same-type numeric value updates ---
$ws.Range("I15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("G16").Value = 1
$ws.Range("N16").Value = -69.230769230769
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -41.666666666666
$ws.Range("I17").Value = 28
$ws.Range("J17").Value = 31
$ws.Range("K17").Value = -9.677419354838
$ws.Range("L17").Value = 55.555555555555
$ws.Range("M17").Value = 133.333333333333
$ws.Range("N17").Value = 100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 280
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -26.923076923076
$ws.Range("N18").Value = -67.241379310344
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 68
$ws.Range("K19").Value = -19.117647058823
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = 37.5
$ws.Range("N19").Value = 61.764705882352
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 0
$ws.Range("N20").Value = -96.491228070175
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -52.941176470588
$ws.Range("F21").Value = 37
$ws.Range("G21").Value = 37
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 114
$ws.Range("J21").Value = 124
$ws.Range("K21").Value = -8.064516129032
$ws.Range("L21").Value = -9.523809523809
$ws.Range("M21").Value = 25.274725274725
$ws.Range("N21").Value = -61.092150170648
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 23
$ws.Range("G24").Value = 31
$ws.Range("H24").Value = -25.806451612903
$ws.Range("I24").Value = 98
$ws.Range("J24").Value = 99
$ws.Range("K24").Value = -1.010101010101
$ws.Range("L24").Value = -26.315789473684
$ws.Range("M24").Value = -9.259259259259
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -40.909090909090
$ws.Range("I25").Value = 50
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = -35.064935064935
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 44.444444444444
$ws.Range("I26").Value = 56
$ws.Range("J26").Value = 32
$ws.Range("K26").Value = 75
$ws.Range("L26").Value = 19.148936170212
$ws.Range("M26").Value = 21.739130434782
$ws.Range("I27").Value = 2
$ws.Range("K27").Value = -60
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 150
